{"js": "// Edit 1: Collapse the three CORE COMPETENCIES detail paragraphs into a\n// single summary paragraph, removing the two now-redundant paragraphs.\n// Edit 2: Append a new \"TECHNICAL SKILLS\" section (heading + three detail\n// paragraphs) at the end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the three CORE COMPETENCIES detail paragraphs by their distinctive\n// leading text so the script does not depend on brittle fixed indices.\nlet pmsIndex = -1; // \"Product Management & Strategy: ...\"\nlet tpdIndex = -1; // \"Technical Product Development: ...\"\nlet pniIndex = -1; // \"Platform & Infrastructure: ...\"\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (pmsIndex === -1 && t.indexOf(\"Product Management & Strategy: Product Conception & Ideation\") === 0) {\n    pmsIndex = i;\n  } else if (tpdIndex === -1 && t.indexOf(\"Technical Product Development: Full-Stack Development\") === 0) {\n    tpdIndex = i;\n  } else if (pniIndex === -1 && t.indexOf(\"Platform & Infrastructure: Multi-tenant Architecture\") === 0) {\n    pniIndex = i;\n  }\n}\n\nif (pmsIndex === -1 || tpdIndex === -1 || pniIndex === -1) {\n  throw new Error(\"Could not locate CORE COMPETENCIES detail paragraphs\");\n}\n\n// Replace the first paragraph's text with the condensed summary line, then\n// delete the other two (now-duplicate) detail paragraphs.\nitems[pmsIndex].insertText(\n  \"Product Management & Strategy \\u2022 Technical Product Development \\u2022 Platform & Infrastructure\",\n  Word.InsertLocation.replace\n);\nitems[tpdIndex].delete();\nitems[pniIndex].delete();\nawait context.sync();\n\n// Append the new TECHNICAL SKILLS section at the very end of the document.\n// Insert the three plain body paragraphs first (they pick up the \"Normal\"\n// style implicitly from the last bulleted paragraph, with no explicit\n// pStyle element, matching the source formatting), then insert the\n// Heading2-styled title immediately before them.\nconst lastParas = body.paragraphs;\nlastParas.load(\"items\");\nawait context.sync();\nconst lastItems = lastParas.items;\nconst lastParagraph = lastItems[lastItems.length - 1];\n\nconst p1 = lastParagraph.insertParagraph(\n  \"PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics\",\n  Word.InsertLocation.after\n);\n\nconst p2 = p1.insertParagraph(\n  \"TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration\",\n  Word.InsertLocation.after\n);\n\nconst p3 = p2.insertParagraph(\n  \"PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training\",\n  Word.InsertLocation.after\n);\n\nconst heading = p1.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.before);\nheading.style = \"Heading2\";\n\nawait context.sync();\n", "ps1": "# Edit 1: Collapse the three CORE COMPETENCIES detail paragraphs into a\n# single summary paragraph, removing the two now-redundant paragraphs.\n# Edit 2: Append a new \"TECHNICAL SKILLS\" section (heading + three detail\n# paragraphs) at the end of the document body.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: condense CORE COMPETENCIES paragraphs -------------------------\n\n$bullet = [char]0x2022\n\n$pmsIndex = -1\n$tpdIndex = -1\n$pniIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($pmsIndex -eq -1 -and $t.StartsWith(\"Product Management & Strategy: Product Conception & Ideation\")) {\n        $pmsIndex = $i\n    } elseif ($tpdIndex -eq -1 -and $t.StartsWith(\"Technical Product Development: Full-Stack Development\")) {\n        $tpdIndex = $i\n    } elseif ($pniIndex -eq -1 -and $t.StartsWith(\"Platform & Infrastructure: Multi-tenant Architecture\")) {\n        $pniIndex = $i\n    }\n}\n\nif ($pmsIndex -eq -1 -or $tpdIndex -eq -1 -or $pniIndex -eq -1) {\n    throw \"Could not locate CORE COMPETENCIES detail paragraphs\"\n}\n\n$pms = $d.Paragraphs.Item($pmsIndex)\n$pms.Range.Text = \"Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure\"\n\n# The other two detail paragraphs immediately follow; deleting one shifts the\n# next one up to the same index, so delete the same index twice.\n$d.Paragraphs.Item($tpdIndex).Range.Delete()\n$d.Paragraphs.Item($tpdIndex).Range.Delete()\n\n# --- Edit 2: append TECHNICAL SKILLS section --------------------------------\n\n$last = $d.Paragraphs.Last\n$last.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics\"\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration\"\n\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Last\n$p3.Range.Text = \"PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training\"\n\n# Insert the heading before the first body paragraph. Re-fetch fresh\n# references afterwards since InsertParagraphBefore() can leave earlier\n# paragraph object handles stale.\n$p3fresh = $d.Paragraphs.Last\n$p2fresh = $p3fresh.Previous()\n$p1fresh = $p2fresh.Previous()\n$p1fresh.Range.InsertParagraphBefore()\n\n$p3fresh2 = $d.Paragraphs.Last\n$p2fresh2 = $p3fresh2.Previous()\n$p1fresh2 = $p2fresh2.Previous()\n$heading = $p1fresh2.Previous()\n$heading.Range.Text = \"TECHNICAL SKILLS\"\n$heading.Style = \"Heading 2\"\n"}
